$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# First-page header (header1.xml): BTec logo, id=1  image2.jpg -> image1.jpg
$h2 = $sec.Headers.Item(2)
if ($h2.Exists -and $h2.Range.InlineShapes.Count -gt 0) {
    $btecLogo = $h2.Range.InlineShapes.Item(1)
    $btecLogo.Name = "image1.jpg"
}

# Default (primary) footer (footer2.xml): Pearson logo, id=2  image1.png -> image2.png
$f1 = $sec.Footers.Item(1)
if ($f1.Exists -and $f1.Range.InlineShapes.Count -gt 0) {
    $pearsonLogoDefault = $f1.Range.InlineShapes.Item(1)
    $pearsonLogoDefault.Name = "image2.png"
}

# First-page footer (footer1.xml): Pearson logo, id=3  image1.png -> image2.png
$f2 = $sec.Footers.Item(2)
if ($f2.Exists -and $f2.Range.InlineShapes.Count -gt 0) {
    $pearsonLogoFirst = $f2.Range.InlineShapes.Item(1)
    $pearsonLogoFirst.Name = "image2.png"
}
